# Bump the "Förändrad" (Changed) date in column C for every data row
# (rows 2 through 347) from 2023-09-09 (serial 45178) to 2023-09-10
# (serial 45179), matching an automatic daily-refresh update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 347

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}
